# metadataframsteg.xlsx - "avtstånd närmsta föryngring klart."
#
# The "GIS-data lyornas avstånd andra fjällrävslyor" task (row 11) is now
# finished: its status moves from "påbörjat" to "klar" (reusing the same
# green "Brödtext" look already used for the other "klar" rows) and its
# comment is replaced with a fresh note about the remaining caveat.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11: "GIS-data lyornas avstånd andra fjällrävslyor" -------------
# Copy the formatting already used by the other "klar" cells (B9) onto
# B11, then set its value - this reuses the existing green/"Brödtext"
# style instead of inventing a new one.
$ws.Range("B9").Copy() | Out-Null
$ws.Range("B11").PasteSpecial(-4122) | Out-Null          # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("B11").Value = "klar"

# Replace the outdated comment with the updated status note.
$ws.Range("C11").Value = "Jag borde ha alla avstånd nu. Det kan hända att jag har något avstånd för mycket ifall det visar sig att en viss lya inte hade kull. Gör om i så fall. "

# --- Window / selection bookkeeping -------------------------------------
# Last selected cell moved to C14.
$ws.Range("C14").Select() | Out-Null

# Remember the (cosmetic) window position at the time of the edit.
$excel.Windows.Item(1).Top = 3200
$excel.Windows.Item(1).Left = 160
